# Update latest output (run 178)
# Refreshes the "Schedule" sheet (adds a new pump-run row, rewrites the others)
# and the "Detailed" sheet (updated prices / pump status / historical-vs-forecast
# classification) with the latest optimisation results.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")

# Insert a new row 4 (pushes old row 4 data down to row 5, duplicating style from row above)
$ws1.Rows.Item(4).Insert()

# Row 2
$ws1.Range("A2").Value = 46068.04166666666
$ws1.Range("B2").Value = 46068.20833333334
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 15.12
$ws1.Range("E2").Value = 428.757069
$ws1.Range("F2").Value = 28.35694900793651

# Row 3
$ws1.Range("A3").Value = 46068.29166666666
$ws1.Range("B3").Value = 46068.77083333334
$ws1.Range("C3").Value = 11.5
$ws1.Range("D3").Value = 43.47
$ws1.Range("E3").Value = 390.83094375
$ws1.Range("F3").Value = 8.990819962042789

# Row 4
$ws1.Range("A4").Value = 46068.97916666666
$ws1.Range("B4").Value = 46069.14583333334
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 15.12
$ws1.Range("E4").Value = 389.6062657499999
$ws1.Range("F4").Value = 25.76761016865079

# Row 5
$ws1.Range("A5").Value = 46069.3125
$ws1.Range("B5").Value = 46069.66666666666
$ws1.Range("C5").Value = 8.5
$ws1.Range("D5").Value = 32.13
$ws1.Range("E5").Value = 41.550171
$ws1.Range("F5").Value = 1.293189262371615

$ws2 = $wb.Worksheets.Item("Detailed")

$ws2.Range("E4").Value = "ON"
$ws2.Range("E5").Value = "ON"
$ws2.Range("E6").Value = "ON"
$ws2.Range("E7").Value = "ON"
$ws2.Range("E8").Value = "ON"
$ws2.Range("E12").Value = "OFF"
$ws2.Range("E13").Value = "OFF"
$ws2.Range("E14").Value = "OFF"
$ws2.Range("E15").Value = "OFF"
$ws2.Range("B37").Value = 12.9589
$ws2.Range("B38").Value = 36.89097
$ws2.Range("B39").Value = 58.05825
$ws2.Range("C40").Value = "historical"
$ws2.Range("B41").Value = 69.09273
$ws2.Range("C41").Value = "historical"
$ws2.Range("B42").Value = 65
$ws2.Range("C42").Value = "historical"
$ws2.Range("B43").Value = 65
$ws2.Range("C43").Value = "historical"
$ws2.Range("B44").Value = 64.89
$ws2.Range("C44").Value = "historical"
$ws2.Range("B45").Value = 60.63248
$ws2.Range("C45").Value = "historical"
$ws2.Range("B46").Value = 57.06007
$ws2.Range("C46").Value = "historical"
$ws2.Range("B47").Value = 60.36903
$ws2.Range("C47").Value = "historical"
$ws2.Range("B48").Value = 50.15088
$ws2.Range("C48").Value = "historical"
$ws2.Range("E48").Value = "OFF"
$ws2.Range("B49").Value = 48.57152
$ws2.Range("B50").Value = 53.39786
$ws2.Range("B52").Value = 36.0601
$ws2.Range("B53").Value = 51.67385
$ws2.Range("B54").Value = 51.54847
$ws2.Range("B55").Value = 50.52186
$ws2.Range("B56").Value = 50.84251
$ws2.Range("B57").Value = 52.21423
$ws2.Range("B58").Value = 54.75563
$ws2.Range("B59").Value = 56.98
$ws2.Range("B60").Value = 57.31
$ws2.Range("B61").Value = 64.89
$ws2.Range("B62").Value = 73.19
$ws2.Range("B64").Value = 57.31
$ws2.Range("B65").Value = 36.0601
$ws2.Range("B66").Value = 4.57358
$ws2.Range("B67").Value = 0.51
$ws2.Range("B68").Value = 0.51
$ws2.Range("B69").Value = 0.51
$ws2.Range("B70").Value = 0.01012
$ws2.Range("B71").Value = 0.36345
$ws2.Range("B72").Value = 0.51
$ws2.Range("B73").Value = 0.46387
$ws2.Range("B74").Value = 0
$ws2.Range("B75").Value = -3.75985
$ws2.Range("B76").Value = -3.75985
$ws2.Range("B77").Value = -2.94057
$ws2.Range("B78").Value = 0.50986
$ws2.Range("B79").Value = 0.98597
$ws2.Range("B80").Value = 0.51
$ws2.Range("B81").Value = 7.55888
$ws2.Range("B82").Value = 31.797
$ws2.Range("B83").Value = 39.69506
$ws2.Range("B84").Value = 52.04194
$ws2.Range("B85").Value = 57.06007
$ws2.Range("B87").Value = 76.07167
$ws2.Range("B88").Value = 90.56927
$ws2.Range("B89").Value = 84.79000000000001
$ws2.Range("B90").Value = 79.92312
$ws2.Range("B92").Value = 72.4654
$ws2.Range("B93").Value = 71.40000000000001
$ws2.Range("B94").Value = 62.65376
$ws2.Range("B95").Value = 57.09
$ws2.Range("B96").Value = 57.09
$ws2.Range("B97").Value = 57.09